$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 454721.38
$ws.Range("I6").Value = 1250022.5
$ws.Range("J6").Value = 263.57144
$ws.Range("K6").Value = 3750067.5
$ws.Range("L6").Value = 790.71432
$ws.Range("M6").Value = -3749955.5
$ws.Range("N6").Value = -1014.71432
$ws.Range("H9").Value = 1819093.5
$ws.Range("I9").Value = 4000528
$ws.Range("J9").Value = 1231.5
$ws.Range("K9").Value = 4000528
$ws.Range("L9").Value = 1231.5
$ws.Range("M9").Value = -4000359
$ws.Range("N9").Value = -1569.5
$ws.Range("H12").Value = 237.14285
$ws.Range("I12").Value = 152.75
$ws.Range("J12").Value = 349.66666
$ws.Range("K12").Value = 152.75
$ws.Range("L12").Value = 349.66666
$ws.Range("M12").Value = 17.25
$ws.Range("N12").Value = -689.66666
$ws.Range("H21").Value = 31000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2936
$ws.Range("H23").Value = 31000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2468
$ws.Range("H42").Value = 1420.5555
$ws.Range("I42").Value = 71.25
$ws.Range("J42").Value = 2500
$ws.Range("K42").Value = 213.75
$ws.Range("L42").Value = 7500
$ws.Range("M42").Value = 16.25
$ws.Range("N42").Value = -7960
$ws.Range("H100").Value = 1929.2858
$ws.Range("I100").Value = 1641
$ws.Range("K100").Value = 1641
$ws.Range("M100").Value = -1100
$ws.Range("H106").Value = 2176.6538
$ws.Range("I106").Value = 2233.0417
$ws.Range("K106").Value = 2233.0417
$ws.Range("M106").Value = -1602.0417
$ws.Range("H138").Value = 6437.3
$ws.Range("I138").Value = 6198
$ws.Range("J138").Value = 6449.8945
$ws.Range("K138").Value = 18594
$ws.Range("L138").Value = 19349.6835
$ws.Range("M138").Value = -13454
$ws.Range("N138").Value = -29629.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5903996
$ws.Range("I32").Value = 5903996
$ws.Range("K32").Value = 5903996
$ws.Range("M32").Value = -5903709
$ws.Range("H61").Value = 28579260
$ws.Range("I61").Value = 4827.8335
$ws.Range("K61").Value = 4827.8335
$ws.Range("M61").Value = -4615.8335
$ws.Range("H74").Value = 51539.24
$ws.Range("I74").Value = 103031.4
$ws.Range("J74").Value = 4728.1816
$ws.Range("K74").Value = 103031.4
$ws.Range("L74").Value = 4728.1816
$ws.Range("M74").Value = -102157.4
$ws.Range("N74").Value = -6476.1816
$ws.Range("H77").Value = 51539.24
$ws.Range("I77").Value = 103031.4
$ws.Range("J77").Value = 4728.1816
$ws.Range("K77").Value = 515157
$ws.Range("L77").Value = 23640.908
$ws.Range("M77").Value = -510789
$ws.Range("N77").Value = -32376.908
$ws.Range("H122").Value = 4245.1113
$ws.Range("I122").Value = 3775.75
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 11327.25
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -8877.25
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 6896
$ws.Range("I132").Value = 3456.7693
$ws.Range("K132").Value = 10370.3079
$ws.Range("M132").Value = -7840.3079
$ws.Range("H136").Value = 28579260
$ws.Range("I136").Value = 4827.8335
$ws.Range("K136").Value = 14483.5005
$ws.Range("M136").Value = -11933.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 10418043
$ws.Range("I64").Value = 22223052
$ws.Range("J64").Value = 1858.0588
$ws.Range("K64").Value = 22223052
$ws.Range("L64").Value = 1858.0588
$ws.Range("M64").Value = -22222827
$ws.Range("N64").Value = -2308.0588
$ws.Range("H67").Value = 10418043
$ws.Range("I67").Value = 22223052
$ws.Range("J67").Value = 1858.0588
$ws.Range("K67").Value = 22223052
$ws.Range("L67").Value = 1858.0588
$ws.Range("M67").Value = -22222272
$ws.Range("N67").Value = -3418.0588
$ws.Range("H134").Value = 5214593
$ws.Range("J134").Value = 9075.68
$ws.Range("L134").Value = 27227.04
$ws.Range("N134").Value = -32297.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18821.646
$ws.Range("I31").Value = 6497.3335
$ws.Range("J31").Value = 25544
$ws.Range("K31").Value = 6497.3335
$ws.Range("L31").Value = 25544
$ws.Range("M31").Value = -6202.3335
$ws.Range("N31").Value = -26134
$ws.Range("H34").Value = 18821.646
$ws.Range("I34").Value = 6497.3335
$ws.Range("J34").Value = 25544
$ws.Range("K34").Value = 6497.3335
$ws.Range("L34").Value = 25544
$ws.Range("M34").Value = -6295.3335
$ws.Range("N34").Value = -25948
$ws.Range("H122").Value = 16669344
$ws.Range("I122").Value = 20835612
$ws.Range("K122").Value = 62506836
$ws.Range("M122").Value = -62504386
$ws.Range("H132").Value = 7043.364
$ws.Range("I132").Value = 4966.864
$ws.Range("K132").Value = 14900.592
$ws.Range("M132").Value = -12370.592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77475.766
$ws.Range("I2").Value = 207.58824
$ws.Range("J2").Value = 223426.78
$ws.Range("K2").Value = 1245.52944
$ws.Range("L2").Value = 1340560.68
$ws.Range("M2").Value = -1132.52944
$ws.Range("N2").Value = -1340786.68
$ws.Range("H34").Value = 3425.3
$ws.Range("J34").Value = 11295.667
$ws.Range("L34").Value = 33887.001
$ws.Range("N34").Value = -34055.001
$ws.Range("H39").Value = 14549.2
$ws.Range("J39").Value = 16999
$ws.Range("L39").Value = 50997
$ws.Range("N39").Value = -51585
$ws.Range("H132").Value = 2619.2917
$ws.Range("I132").Value = 2048.3333
$ws.Range("J132").Value = 4332.1665
$ws.Range("K132").Value = 18434.9997
$ws.Range("L132").Value = 38989.4985
$ws.Range("M132").Value = -15904.9997
$ws.Range("N132").Value = -44049.4985
$ws.Range("H140").Value = 225103.44
$ws.Range("J140").Value = 3704.7144
$ws.Range("L140").Value = 11114.1432
$ws.Range("N140").Value = -21474.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1603
$ws.Range("I97").Value = 1482.0526
$ws.Range("K97").Value = 1482.0526
$ws.Range("M97").Value = -986.0526
$ws.Range("H132").Value = 4774.892
$ws.Range("I132").Value = 2867.2
$ws.Range("K132").Value = 8601.599999999999
$ws.Range("M132").Value = -6071.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 772.5333000000001
$ws.Range("I55").Value = 602.3333
$ws.Range("J55").Value = 886
$ws.Range("K55").Value = 602.3333
$ws.Range("L55").Value = 886
$ws.Range("M55").Value = -429.3333
$ws.Range("N55").Value = -1232
$ws.Range("H68").Value = 6267.3335
$ws.Range("J68").Value = 11000
$ws.Range("L68").Value = 11000
$ws.Range("N68").Value = -12498
$ws.Range("H71").Value = 6267.3335
$ws.Range("J71").Value = 11000
$ws.Range("L71").Value = 55000
$ws.Range("N71").Value = -62488
$ws.Range("H82").Value = 502222
$ws.Range("I82").Value = 1000000
$ws.Range("K82").Value = 1000000
$ws.Range("M82").Value = -999639
$ws.Range("H85").Value = 502222
$ws.Range("I85").Value = 1000000
$ws.Range("K85").Value = 1000000
$ws.Range("M85").Value = -998752
$ws.Range("H122").Value = 6148.0464
$ws.Range("J122").Value = 6687
$ws.Range("L122").Value = 20061
$ws.Range("N122").Value = -24961
$ws.Range("H132").Value = 20843258
$ws.Range("I132").Value = 38466708
$ws.Range("J132").Value = 15545.272
$ws.Range("K132").Value = 115400124
$ws.Range("L132").Value = 46635.81600000001
$ws.Range("M132").Value = -115397594
$ws.Range("N132").Value = -51695.81600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 20000
$ws.Range("I76").Value = 20000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -19685
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 20000
$ws.Range("I79").Value = 20000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -18908
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 9089.531000000001
$ws.Range("I113").Value = 13151
$ws.Range("J113").Value = 2320.4167
$ws.Range("K113").Value = 39453
$ws.Range("L113").Value = 6961.250100000001
$ws.Range("M113").Value = -37283
$ws.Range("N113").Value = -11301.2501
$ws.Range("H132").Value = 13026.742
$ws.Range("I132").Value = 10333.16
$ws.Range("J132").Value = 24250
$ws.Range("K132").Value = 30999.48
$ws.Range("L132").Value = 72750
$ws.Range("M132").Value = -28469.48
$ws.Range("N132").Value = -77810
$ws.Range("H135").Value = 78532
$ws.Range("J135").Value = 78532
$ws.Range("L135").Value = 78532
$ws.Range("N135").Value = -88672
$ws.Range("H136").Value = 24052086
$ws.Range("I136").Value = 45456616
$ws.Range("K136").Value = 136369848
$ws.Range("M136").Value = -136367298
